$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Conditions of Release:", $true, $false, $false, $false, $false,
               $true, 1, $false, "Requirements for Release:", 2)
